# Insert a new data row at row 194, pushing the existing rows 194-266 down
# to 195-267 (dimension grows from A1:R266 to A1:R267).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 194; this shifts rows 194:266 down
# to 195:267 and carries their formatting/styles along automatically.
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row 194 with its values.
$ws.Cells.Item(194, 1).Value = 6
$ws.Cells.Item(194, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(194, 3).Value = "Metropolitana"
$ws.Cells.Item(194, 4).Value = 44900
$ws.Cells.Item(194, 5).Value = 13
$ws.Cells.Item(194, 6).Value = 100112029
$ws.Cells.Item(194, 7).Value = "Orégano"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 34
$ws.Cells.Item(194, 11).Value = 16000
$ws.Cells.Item(194, 12).Value = 17000
$ws.Cells.Item(194, 13).Value = 16441
$ws.Cells.Item(194, 14).Value = "`$/docena de atados"
$ws.Cells.Item(194, 15).Value = "Región Metropolitana"
$ws.Cells.Item(194, 16).Value = 5480
$ws.Cells.Item(194, 17).Value = 3
$ws.Cells.Item(194, 18).Value = "Hortaliza"

# Make sure the new row's date cell keeps the same custom date number format
# as the rest of column D.
$ws.Cells.Item(194, 4).NumberFormat = $ws.Cells.Item(195, 4).NumberFormat
